$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Swap country names for rows 217 and 218 (Comoras <-> San Pedro y Miquelon)
$ws.Range("A217").Value = "San Pedro y Miquelon"
$ws.Range("A218").Value = "Comoras"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1152328
$ws.Range("C4").Value = 21298
$ws.Range("D4").Value = 170179
$ws.Range("E4").Value = 915228
$ws.Range("G4").Value = 1168
$ws.Range("H4").Value = 66921

# Row 17: Peru
$ws.Range("B17").Value = 42534
$ws.Range("C17").Value = 2075
$ws.Range("D17").Value = 12434
$ws.Range("E17").Value = 28900
$ws.Range("F17").Value = 671
$ws.Range("G17").Value = 76
$ws.Range("H17").Value = 1200

# Row 30: Israel
$ws.Range("B30").Value = 16185
$ws.Range("C30").Value = 84
$ws.Range("D30").Value = 9593
$ws.Range("E30").Value = 6363
$ws.Range("F30").Value = 105
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 229

# Row 45: Noruega
$ws.Range("E45").Value = 7558
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 211

# Row 124: Venezuela
$ws.Range("B124").Value = 345
$ws.Range("C124").Value = 10
$ws.Range("E124").Value = 187
